# feat: add 2022-Q4 data
#
# Before:  总计 | 2022-Q3 | 2020-Q4
# After:   总计 | 2022-Q4 | 2022-Q3 | 2020-Q4
#
# The old "2022-Q3" sheet is duplicated (the duplicate keeps the old,
# unmodified numbers and is renamed back to "2022-Q3"); the original
# sheet is renamed to "2022-Q4" and its figures are refreshed. The
# "总计" (totals) sheet gets a new row for the 2022-Q4 quarter, pushing
# the existing rows down.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# --- 1. Duplicate the "2022-Q3" sheet right after itself -----------------
# The copy preserves the current (pre-edit) figures; it will become the
# "new" 2022-Q3 sheet, while the original sheet is turned into 2022-Q4.
$q3.Copy($null, $q3)
$dup = $wb.Worksheets.Item("2022-Q3 (2)")

# Rename: original -> 2022-Q4 (updated numbers), duplicate -> 2022-Q3 (kept numbers)
$dup.Name = "2022-Q3-staging"
$q3.Name = "2022-Q4"
$dup.Name = "2022-Q3"

$q4 = $wb.Worksheets.Item("2022-Q4")

# --- 2. Refresh the figures on the new 2022-Q4 sheet ----------------------
$q4.Range("D2").Value = "'3.92"
$q4.Range("E2").Value = "'94.38"
$q4.Range("F2").Value = "'1.69"
$q4.Range("G2").Value = "'0.0662"
$q4.Range("D2:G2").ClearFormats()
$q4.Range("H2").Value = 8

# --- 3. Insert the new 2022-Q4 row into the "总计" summary sheet ---------
# Copy row 2's formatting down into the new row 3 first, then shift the
# two existing data rows down by one and write the new row 2.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.07000000000000001

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.06

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2020-Q4"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01
